$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productdata")

$ws.Range("C4").Value = 837
$ws.Range("C5").Value = 837
$ws.Range("C6").Value = 837
$ws.Range("C7").Value = 1013
$ws.Range("C8").Value = 181
$ws.Range("C9").Value = 301
